$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07966813408545038
$ws.Range("C2").Value = 0.2904926044902277
$ws.Range("D2").Value = 0.1997595342671157
$ws.Range("E2").Value = 0.4469446657776729
$ws.Range("F2").Value = 0.4440776115544651

$ws.Range("B3").Value = -0.09581970561355245
$ws.Range("C3").Value = 0.4927854440499702
$ws.Range("D3").Value = 0.4359434021580401
$ws.Range("E3").Value = 0.6602601018977597
$ws.Range("F3").Value = 0.6597705858081699

$ws.Range("B4").Value = 0.003072972994658064
$ws.Range("C4").Value = 0.4851683705022662
$ws.Range("D4").Value = 0.4107751511171511
$ws.Range("E4").Value = 0.6409174292505635
$ws.Range("F4").Value = 0.6474169302588882

$ws.Range("B5").Value = -0.1154007161324104
$ws.Range("C5").Value = 0.479353824850625
$ws.Range("D5").Value = 0.4268500492957792
$ws.Range("E5").Value = 0.6533376227462944
$ws.Range("F5").Value = 0.6497291659572375
$ws.Range("G5").Value = 49

$ws.Range("B6").Value = 0.002725215881915507
$ws.Range("C6").Value = 0.4511380015418734
$ws.Range("D6").Value = 0.4051435168807476
$ws.Range("E6").Value = 0.6365088505910563
$ws.Range("F6").Value = 0.6432386857841645
$ws.Range("G6").Value = 48

$ws.Range("B7").Value = -0.04206088621362827
$ws.Range("C7").Value = 0.4214914888921343
$ws.Range("D7").Value = 0.4075230093214146
$ws.Range("E7").Value = 0.6383752887772322
$ws.Range("F7").Value = 0.6453151363098175
$ws.Range("G7").Value = 39

$ws.Range("B8").Value = -0.03202104169280565
$ws.Range("C8").Value = 0.4420571789576535
$ws.Range("D8").Value = 0.4570240603574153
$ws.Range("E8").Value = 0.6760355466670487
$ws.Range("F8").Value = 0.6843412911651034
$ws.Range("G8").Value = 38

$ws.Range("B9").Value = -0.1117388779456075
$ws.Range("C9").Value = 0.4582166133627955
$ws.Range("D9").Value = 0.6029226805016665
$ws.Range("E9").Value = 0.7764809595229406
$ws.Range("F9").Value = 0.7873747258072102
$ws.Range("G9").Value = 21

$ws.Range("B10").Value = 0.01972275674151302
$ws.Range("C10").Value = 0.3574201004013381
$ws.Range("D10").Value = 0.3132652572947969
$ws.Range("E10").Value = 0.5597010427851612
$ws.Range("F10").Value = 0.5804684966114319
$ws.Range("G10").Value = 14

$ws.Range("B11").Value = 0.2039999999999907
$ws.Range("C11").Value = 0.4360000000000014
$ws.Range("D11").Value = 0.4379999999999994
$ws.Range("E11").Value = 0.6618156843109715
$ws.Range("F11").Value = 0.7039034024637215
